# Update the "dSF" (column F) values on Sheet1 to reflect the re-pulled /
# recalculated data described in the commit message ("repull data, push
# all data, mean calculation").
#
# Only column F changes; every other column / row is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of row number -> new dSF (column F) value
$updates = @{
    2  = -8
    3  = -5
    4  = -3
    5  = -16
    9  = -10
    10 = -7
    12 = -7
    14 = -5
    15 = -5
    17 = -5
    19 = -5
    24 = -1
    27 = -6
    29 = 10
    32 = 1
    33 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
